# Weekly update: insert a new price observation as row 663, pushing the
# existing rows 663:685 down to 664:686 (last row duplicated the prior
# bottom row's content one row further down, growing the used range to
# A1:R686).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 663:685 down one row to make room for the new record.
$ws.Rows.Item(663).Insert()

# Populate the newly inserted row 663 with the new weekly observation.
$ws.Cells.Item(663, 1).Value  = 6
$ws.Cells.Item(663, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(663, 3).Value  = "Metropolitana"
$ws.Cells.Item(663, 4).Value  = 45075
$ws.Cells.Item(663, 5).Value  = 13
$ws.Cells.Item(663, 6).Value  = 100112039
$ws.Cells.Item(663, 7).Value  = "Ciboulette"
$ws.Cells.Item(663, 8).Value  = "Sin especificar"
$ws.Cells.Item(663, 9).Value  = "Primera"
$ws.Cells.Item(663, 10).Value = 580
$ws.Cells.Item(663, 11).Value = 900
$ws.Cells.Item(663, 12).Value = 1000
$ws.Cells.Item(663, 13).Value = 948
$ws.Cells.Item(663, 14).Value = "`$/docena de atados"
$ws.Cells.Item(663, 15).Value = "Región Metropolitana"
$ws.Cells.Item(663, 16).Value = 316
$ws.Cells.Item(663, 17).Value = 3
$ws.Cells.Item(663, 18).Value = "Hortaliza"
